$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new row of pledge data (string literals must be entered in this
# order so the shared-strings table indices come out as 18 then 19)
$ws.Range("B7").Value = "(t+auto) GO Bedrock 2018"
$ws.Range("A7").Value = "mpcrds+auto+fredflintstone@gmail.com"
$ws.Range("C7").Value = 1000
$ws.Range("D7").Value = 0
$ws.Range("E7").Value = 43101

# Match date formatting used by the other date cells in column E
$ws.Range("E6").Copy()
$ws.Range("E7").PasteSpecial(-4122)

# Reflect the UI selection state recorded in the saved file
$ws.Range("D15").Select()
